$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B, shifting all existing price-history
# columns (old B..U) one slot to the right (new C..V). This also carries the
# per-column width/style formatting of column B to the new columns.
$ws.Columns("B").Insert()

# Make sure the freshly inserted column keeps the same width as its neighbours.
$ws.Columns("B").ColumnWidth = $ws.Columns("C").ColumnWidth

# New timestamp header for this scrape.
$ws.Range("B1").Value = "2025-12-22 10:32"

# New price observations for this scrape time (row 1 is the header, handled above).
$ws.Range("B2").Value = 929
$ws.Range("B4").Value = 299
$ws.Range("B5").Value = 569
$ws.Range("B6").Value = 499
$ws.Range("B7").Value = 569
$ws.Range("B8").Value = 929
$ws.Range("B9").Value = 299
$ws.Range("B10").Value = 299
$ws.Range("B11").Value = 2997
$ws.Range("B12").Value = 569
$ws.Range("B13").Value = 569
$ws.Range("B14").Value = 499
$ws.Range("B15").Value = 499
$ws.Range("B16").Value = 299
$ws.Range("B17").Value = 929
$ws.Range("B18").Value = 499
$ws.Range("B19").Value = 1299
$ws.Range("B20").Value = 929
$ws.Range("B21").Value = 499
$ws.Range("B22").Value = 299
$ws.Range("B23").Value = 1299
$ws.Range("B24").Value = 929
$ws.Range("B25").Value = 929
$ws.Range("B26").Value = 1299
# B3 has no price for this scrape (out of stock) and is intentionally left blank.
